# Insert a new data row before current row 452, shifting rows 452:546 down to 453:547,
# and populate the newly inserted row with the "Provincia de Chacabuco" / "Primera" entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(452).Insert()

$row = 452

$ws.Cells.Item($row, 1).Value2 = 9
$ws.Cells.Item($row, 2).Value2 = "Vega Central Mapocho de Santiago"
$ws.Cells.Item($row, 3).Value2 = "Metropolitana"
$ws.Cells.Item($row, 4).Value2 = 44889
$ws.Cells.Item($row, 5).Value2 = 13
$ws.Cells.Item($row, 6).Value2 = 100112013
$ws.Cells.Item($row, 7).Value2 = "Alcachofa"
$ws.Cells.Item($row, 8).Value2 = "Española"
$ws.Cells.Item($row, 9).Value2 = "Primera"
$ws.Cells.Item($row, 10).Value2 = 150
$ws.Cells.Item($row, 11).Value2 = 8000
$ws.Cells.Item($row, 12).Value2 = 8000
$ws.Cells.Item($row, 13).Value2 = 8000
$ws.Cells.Item($row, 14).Value2 = "$/caja 30 unidades"
$ws.Cells.Item($row, 15).Value2 = "Provincia de Chacabuco"
$ws.Cells.Item($row, 16).Value2 = 267
$ws.Cells.Item($row, 17).Value2 = 30
$ws.Cells.Item($row, 18).Value2 = "Hortaliza"
